$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5991200276118605
$ws.Range("C2").Value = 0.1697399270884432
$ws.Range("D2").Value = 0.07956700250065296
$ws.Range("E2").Value = 0.1329857292504002
$ws.Range("G2").Value = 0.3087939921585274
$ws.Range("H2").Value = 0.4533695470216799
$ws.Range("I2").Value = 0.3157090693811675
$ws.Range("M2").Value = 0.3011909661693366
$ws.Range("N2").Value = 0.9113989737215178
$ws.Range("O2").Value = 1.452456590346543
$ws.Range("B3").Value = 0.5236023706505648
$ws.Range("C3").Value = 0.1503981895707227
$ws.Range("D3").Value = 0.07203563985159178
$ws.Range("E3").Value = 0.1264605457262675
$ws.Range("G3").Value = 0.3021107930165527
$ws.Range("H3").Value = 0.4541211505721634
$ws.Range("I3").Value = 0.3186110459157732
$ws.Range("M3").Value = 0.2676756553880324
$ws.Range("N3").Value = 0.9187644222937976
$ws.Range("O3").Value = 1.439573425276791
$ws.Range("B4").Value = 0.4771125849678128
$ws.Range("C4").Value = 0.1384557726621836
$ws.Range("D4").Value = 0.06744487001265043
$ws.Range("E4").Value = 0.1225733953228669
$ws.Range("G4").Value = 0.2982940009326853
$ws.Range("H4").Value = 0.4548499962488961
$ws.Range("I4").Value = 0.3206390374463517
$ws.Range("M4").Value = 0.2471461241541633
$ws.Range("N4").Value = 0.9236801844738025
$ws.Range("O4").Value = 1.432820275224572
$ws.Range("B5").Value = 0.4581381430849092
$ws.Range("C5").Value = 0.1335726409637061
$ws.Range("D5").Value = 0.06558252851750979
$ws.Range("E5").Value = 0.1210191052567922
$ws.Range("G5").Value = 0.2968104808139316
$ws.Range("H5").Value = 0.4552142098076359
$ws.Range("I5").Value = 0.3215272345280695
$ws.Range("M5").Value = 0.2387925623148703
$ws.Range("N5").Value = 0.9257823688332749
$ws.Range("O5").Value = 1.430358710282547
$ws.Range("B6").Value = 0.4549857038853702
$ws.Range("C6").Value = 0.1327608107579579
$ws.Range("D6").Value = 0.0652737977423925
$ws.Range("E6").Value = 0.1207628063175861
$ws.Range("G6").Value = 0.2965684760571747
$ws.Range("H6").Value = 0.4552787453282718
$ws.Range("I6").Value = 0.3216784462735909
$ws.Range("M6").Value = 0.2374062099950649
$ws.Range("N6").Value = 0.9261374144483199
$ws.Range("O6").Value = 1.429967491683641
$ws.Range("B7").Value = 0.4768568068860191
$ws.Range("C7").Value = 0.1383899835495299
$ws.Range("D7").Value = 0.06741971966449967
$ws.Range("E7").Value = 0.1225523134801918
$ws.Range("G7").Value = 0.2982737030529989
$ws.Range("H7").Value = 0.4548546360938985
$ws.Range("I7").Value = 0.3206507660302655
$ws.Range("M7").Value = 0.2470334148465696
$ws.Range("N7").Value = 0.9237081344535767
$ws.Range("O7").Value = 1.432785902659248
$ws.Range("B8").Value = 0.5731074284447857
$ws.Range("C8").Value = 0.1630848298498222
$ws.Range("D8").Value = 0.07696323364339719
$ws.Range("E8").Value = 0.1307109224938827
$ws.Range("G8").Value = 0.3064299667847621
$ws.Range("H8").Value = 0.4535731927511506
$ws.Range("I8").Value = 0.3166585028950522
$ws.Range("M8").Value = 0.2896247451436409
$ws.Range("N8").Value = 0.9138569766462226
$ws.Range("O8").Value = 1.447773947952811
$ws.Range("B9").Value = 0.7608496201409594
$ws.Range("C9").Value = 0.2109759311953781
$ws.Range("D9").Value = 0.09594469734088307
$ws.Range("E9").Value = 0.1476691886578436
$ws.Range("G9").Value = 0.3247116763272118
$ws.Range("H9").Value = 0.4531833944583212
$ws.Range("I9").Value = 0.3107887422026891
$ws.Range("M9").Value = 0.3735364594047041
$ws.Range("N9").Value = 0.8976575044168982
$ws.Range("O9").Value = 1.486377124928651
$ws.Range("B10").Value = 0.8981307944973764
$ws.Range("C10").Value = 0.2458280698084536
$ws.Range("D10").Value = 0.1100553780321718
$ws.Range("E10").Value = 0.1607316634992983
$ws.Range("G10").Value = 0.3395565554569373
$ws.Range("H10").Value = 0.4541945252757813
$ws.Range("I10").Value = 0.3076787927742757
$ws.Range("M10").Value = 0.43543258628263
$ws.Range("N10").Value = 0.8876540957502854
$ws.Range("O10").Value = 1.520401820139
$ws.Range("B11").Value = 0.9604340237167435
$ws.Range("C11").Value = 0.2616095106208434
$ws.Range("D11").Value = 0.1165111183398153
$ws.Range("E11").Value = 0.1668091599760047
$ws.Range("G11").Value = 0.346620921555413
$ws.Range("H11").Value = 0.4549370065847
$ws.Range("I11").Value = 0.3065269134972937
$ws.Range("M11").Value = 0.463646354398989
$ws.Range("N11").Value = 0.8835150007506272
$ws.Range("O11").Value = 1.537120671741974
$ws.Range("B12").Value = 0.9840045959477379
$ws.Range("C12").Value = 0.267574854102719
$ws.Range("D12").Value = 0.1189610334463254
$ws.Range("E12").Value = 0.16913030599396
$ws.Range("G12").Value = 0.3493410791451907
$ws.Range("H12").Value = 0.4552588370278414
$ws.Range("I12").Value = 0.3061286679005377
$ws.Range("M12").Value = 0.4743383890284321
$ws.Range("N12").Value = 0.8820067778865095
$ws.Range("O12").Value = 1.543630825923259
$ws.Range("B13").Value = 0.9789292630291015
$ws.Range("C13").Value = 0.266290592579395
$ws.Range("D13").Value = 0.1184331664032214
$ws.Range("E13").Value = 0.168629523828244
$ws.Range("G13").Value = 0.3487532374029314
$ws.Range("H13").Value = 0.4551877156245041
$ws.Range("I13").Value = 0.3062127471957616
$ws.Range("M13").Value = 0.4720353095700176
$ws.Range("N13").Value = 0.8823289695297944
$ws.Range("O13").Value = 1.542220772344336
$ws.Range("B14").Value = 0.9623736429442715
$ws.Range("C14").Value = 0.2621004998214573
$ws.Range("D14").Value = 0.1167125688667454
$ws.Range("E14").Value = 0.166999725074966
$ws.Range("G14").Value = 0.3468438064788586
$ws.Range("H14").Value = 0.4549626684525379
$ws.Range("I14").Value = 0.3064933882483345
$ws.Range("M14").Value = 0.464525832846121
$ws.Range("N14").Value = 0.8833897326918105
$ws.Range("O14").Value = 1.537652673110586
$ws.Range("B15").Value = 0.9522299062742263
$ws.Range("C15").Value = 0.2595325401387925
$ws.Range("D15").Value = 0.1156593386019864
$ws.Range("E15").Value = 0.1660040038013122
$ws.Range("G15").Value = 0.3456800976290282
$ws.Range("H15").Value = 0.4548301181538363
$ws.Range("I15").Value = 0.3066702350909836
$ws.Range("M15").Value = 0.4599271127780042
$ws.Range("N15").Value = 0.8840471852059153
$ws.Range("O15").Value = 1.534877922372232
$ws.Range("B16").Value = 0.8940560764587531
$ws.Range("C16").Value = 0.2447952255695327
$ws.Range("D16").Value = 0.1096342179542518
$ws.Range("E16").Value = 0.160337230894541
$ws.Range("G16").Value = 0.339101170298008
$ws.Range("H16").Value = 0.454151691007084
$ws.Range("I16").Value = 0.3077593727342816
$ws.Range("M16").Value = 0.4335898876055353
$ws.Range("N16").Value = 0.8879328748131599
$ws.Range("O16").Value = 1.519334228992648
$ws.Range("B17").Value = 0.8583298838755695
$ws.Range("C17").Value = 0.2357354908196214
$ws.Range("D17").Value = 0.1059473969196176
$ws.Range("E17").Value = 0.1568956993794117
$ws.Range("G17").Value = 0.3351451446818743
$ws.Range("H17").Value = 0.4538078874580833
$ws.Range("I17").Value = 0.3084949517678446
$ws.Range("M17").Value = 0.417447364458198
$ws.Range("N17").Value = 0.8904220065045791
$ws.Range("O17").Value = 1.510116954668916
$ws.Range("B18").Value = 0.8377673906026075
$ws.Range("C18").Value = 0.2305177176087909
$ws.Range("D18").Value = 0.1038302922053873
$ws.Range("E18").Value = 0.1549289412923045
$ws.Range("G18").Value = 0.3328990340635301
$ws.Range("H18").Value = 0.4536367295645931
$ws.Range("I18").Value = 0.3089427694324343
$ws.Range("M18").Value = 0.4081679741074424
$ws.Range("N18").Value = 0.8918924224888016
$ws.Range("O18").Value = 1.504932177641166
$ws.Range("B19").Value = 0.8308029623369748
$ws.Range("C19").Value = 0.2287498998726392
$ws.Range("D19").Value = 0.1031140711249066
$ws.Range("E19").Value = 0.1542652071443982
$ws.Range("G19").Value = 0.332143562274851
$ws.Range("H19").Value = 0.4535833437632135
$ws.Range("I19").Value = 0.3090986360331343
$ws.Range("M19").Value = 0.4050270556664088
$ws.Range("N19").Value = 0.8923969329752737
$ws.Range("O19").Value = 1.503196734056701
$ws.Range("B20").Value = 0.8621344281321512
$ws.Range("C20").Value = 0.23670062618163
$ws.Range("D20").Value = 0.1063395076020157
$ws.Range("E20").Value = 0.1572607378604403
$ws.Range("G20").Value = 0.3355632369782739
$ws.Range("H20").Value = 0.4538417337708012
$ws.Range("I20").Value = 0.3084140873730234
$ws.Range("M20").Value = 0.4191652093150964
$ws.Range("N20").Value = 0.8901530255333299
$ws.Range("O20").Value = 1.511086060751637
$ws.Range("B21").Value = 0.9672370450181802
$ws.Range("C21").Value = 0.2633315252756177
$ws.Range("D21").Value = 0.1172178072597205
$ws.Range("E21").Value = 0.1674778988357346
$ws.Range("G21").Value = 0.3474034282887857
$ws.Range("H21").Value = 0.4550276661704373
$ws.Range("I21").Value = 0.306409926142635
$ws.Range("M21").Value = 0.4667313292729318
$ws.Range("N21").Value = 0.8830765552518329
$ws.Range("O21").Value = 1.538989568982089
$ws.Range("B22").Value = 1.035797086505397
$ws.Range("C22").Value = 0.2806736067187217
$ws.Range("D22").Value = 0.1243580942742426
$ws.Range("E22").Value = 0.1742705324900129
$ws.Range("G22").Value = 0.3554042999537472
$ws.Range("H22").Value = 0.4560398113188882
$ws.Range("I22").Value = 0.3053213383171354
$ws.Range("M22").Value = 0.4978656684050122
$ws.Range("N22").Value = 0.8787965153996637
$ws.Range("O22").Value = 1.558270288398006
$ws.Range("B23").Value = 0.9992176532288113
$ws.Range("C23").Value = 0.27142363640138
$ws.Range("D23").Value = 0.1205443854717885
$ws.Range("E23").Value = 0.1706345516812107
$ws.Range("G23").Value = 0.3511099677255345
$ws.Range("H23").Value = 0.4554779040001051
$ws.Range("I23").Value = 0.3058820450799118
$ws.Range("M23").Value = 0.4812444031505976
$ws.Range("N23").Value = 0.8810493027665629
$ws.Range("O23").Value = 1.547884048642374
$ws.Range("B24").Value = 0.8604144653156141
$ws.Range("C24").Value = 0.2362643171579748
$ws.Range("D24").Value = 0.1061622265822422
$ws.Range("E24").Value = 0.1570956671543016
$ws.Range("G24").Value = 0.335374129439387
$ws.Range("H24").Value = 0.4538263493104751
$ws.Range("I24").Value = 0.3084505685412218
$ws.Range("M24").Value = 0.4183885680354962
$ws.Range("N24").Value = 0.8902745091573081
$ws.Range("O24").Value = 1.51064757172449
$ws.Range("B25").Value = 0.7101721559600946
$ws.Range("C25").Value = 0.1980781703783521
$ws.Range("D25").Value = 0.09078090627765789
$ws.Range("E25").Value = 0.1429769281099809
$ws.Range("G25").Value = 0.3195191595788032
$ws.Range("H25").Value = 0.4530612203367212
$ws.Range("I25").Value = 0.3121660819163523
$ws.Range("M25").Value = 0.3507934831627253
$ws.Range("N25").Value = 0.8823289695297944
$ws.Range("O25").Value = 1.542220772344336
